$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.749.63"
$ws.Range("E2").Value = "  +8.52%  "
$ws.Range("D3").Value = "1.950.44"
$ws.Range("E3").Value = "  +6.77%  "
$ws.Range("D4").Value = "'0.9952"
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'342.91"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("D6").Value = "'0.9956"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "'0.4804"
$ws.Range("E7").Value = "  +5.06%  "
$ws.Range("D8").Value = "'0.4139"
$ws.Range("E8").Value = "  +8.58%  "
$ws.Range("D9").Value = "'48.82"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("D10").Value = "'0.08286"
$ws.Range("E10").Value = "  +5.22%  "
$ws.Range("D11").Value = "'1.051"
$ws.Range("E11").Value = "  +9.00%  "
$ws.Range("D12").Value = "'22.78"
$ws.Range("E12").Value = "  +8.56%  "
$ws.Range("D13").Value = "1.932.38"
$ws.Range("E13").Value = "  +5.83%  "
$ws.Range("D14").Value = "'6.163"
$ws.Range("E14").Value = "  +4.96%  "
$ws.Range("D15").Value = "'7.451"
$ws.Range("E15").Value = "  +5.37%  "
$ws.Range("D16").Value = "'93.15"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "'0.9960"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'0.00001070"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").Value = "'0.06676"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'18.14"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("D21").Value = "'0.9950"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "29.694.78"
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("D23").Value = "'5.636"
$ws.Range("E23").Value = "  +5.96%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("E24").Value = "  +4.78%  "
$ws.Range("D25").Value = "'2.261"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").Value = "2.162.11"
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("D27").Value = "'161.39"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("D28").Value = "'20.27"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "'2.209"
$ws.Range("E29").Value = "  +7.02%  "
$ws.Range("D30").Value = "'5.644"
$ws.Range("E30").Value = "  +6.75%  "
$ws.Range("D31").Value = "'122.75"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("E32").Value = "  +9.52%  "
$ws.Range("D33").Value = "'0.09667"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("D34").Value = "'1.484"
$ws.Range("E34").Value = "  +12.13%  "
$ws.Range("D35").Value = "'3.678"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "'5.512"
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("D37").Value = "'0.02315"
$ws.Range("E37").Value = "  +6.46%  "
$ws.Range("D38").Value = "'0.06247"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("D39").Value = "'8.703"
$ws.Range("E39").Value = "  +7.22%  "
$ws.Range("D40").Value = "'1.206"
$ws.Range("E40").Value = "  +5.65%  "
$ws.Range("D41").Value = "'0.6129"
$ws.Range("E41").Value = "  +6.34%  "
$ws.Range("D42").Value = "'10.76"
$ws.Range("E42").Value = "  +8.06%  "
$ws.Range("D43").Value = "'0.1916"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("D44").Value = "'0.9948"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "'1.298"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D48").Value = "'2.329"
$ws.Range("E48").Value = "  +28.55%  "
$ws.Range("D49").Value = "'2.010"
$ws.Range("E49").Value = "  +7.51%  "
$ws.Range("D50").Value = "'0.07243"
$ws.Range("E50").Value = "  +10.08%  "
$ws.Range("D51").Value = "'113.86"
$ws.Range("E51").Value = "  +3.31%  "

# Row 46/47 swap: EnergySwap now row 46, Decentraland now row 47
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.59"
$ws.Range("E46").Value = "  +5.11%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5739"
$ws.Range("E47").Value = "  +5.78%  "
